$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.762.87"
$ws.Range("E2").Value = "  +2.93%  "

$ws.Range("D3").Value = "1.863.92"
$ws.Range("E3").Value = "  +2.68%  "

$ws.Range("E4").Value = "  +3.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.07"
$ws.Range("E5").Value = "  +3.43%  "

$ws.Range("E6").Value = "  +2.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4426"
$ws.Range("E7").Value = "  +2.98%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3796"
$ws.Range("E8").Value = "  +3.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07460"
$ws.Range("E9").Value = "  +3.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8847"
$ws.Range("E10").Value = "  +2.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.68"
$ws.Range("E11").Value = "  +2.34%  "

$ws.Range("D12").Value = "1.872.39"
$ws.Range("E12").Value = "  -15.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.557"
$ws.Range("E13").Value = "  +2.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.764"

$ws.Range("E15").Value = "  +4.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.77"
$ws.Range("E16").Value = "  +3.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.042"
$ws.Range("E17").Value = "  +2.87%  "

$ws.Range("E18").Value = "  +2.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.56"
$ws.Range("E20").Value = "  +2.40%  "

$ws.Range("D21").Value = "27.774.35"
$ws.Range("E21").Value = "  +2.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.313"
$ws.Range("E22").Value = "  +2.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.34"
$ws.Range("E23").Value = "  +3.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.006"
$ws.Range("E24").Value = "  +7.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.78"
$ws.Range("E25").Value = "  +3.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.86"
$ws.Range("E26").Value = "  +3.14%  "

$ws.Range("E27").Value = "  +2.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.982"
$ws.Range("E28").Value = "  +4.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "118.07"
$ws.Range("E29").Value = "  +3.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09065"
$ws.Range("E30").Value = "  +1.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7770"
$ws.Range("E31").Value = "  +3.64%  "

$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.215"
$ws.Range("E32").Value = "  +2.21%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.094"
$ws.Range("E33").Value = "  +10.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.572"
$ws.Range("E34").Value = "  +3.42%  "

$ws.Range("E35").Value = "  +3.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.151"
$ws.Range("E36").Value = "  +1.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01994"
$ws.Range("E37").Value = "  +3.63%  "

$ws.Range("E38").Value = "  +2.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.860"
$ws.Range("E39").Value = "  +3.98%  "

$ws.Range("E40").Value = "  +1.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1693"
$ws.Range("E41").Value = "  +2.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.888"
$ws.Range("E42").Value = "  +6.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.668"
$ws.Range("E43").Value = "  +4.17%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.74"
$ws.Range("E44").Value = "  +3.60%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "109.80"
$ws.Range("E45").Value = "  +2.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.720"
$ws.Range("E46").Value = "  +4.59%  "

$ws.Range("E47").Value = "  +3.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06477"
$ws.Range("E48").Value = "  +4.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.918"
$ws.Range("E49").Value = "  +3.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.85"
$ws.Range("E50").Value = "  +1.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.60"
$ws.Range("E51").Value = "  +2.68%  "
